# Generate Report for Handback
#
# The handback pipeline has produced fresh timestamps for the file
# "f8ddade8-8c09-481e-95e5-1efd3d21c69f" (row 3 on every sheet): Excel
# re-generated the HO xliff on the Overview sheet and stamped the
# corresponding handoff/handback datetimes on the per-locale sheets.
#
# Overview!G3            "Latest HO Xliff Generate Date"  -> 2016-09-02 04:51:47
# zh-cn!H3                "Correspond Handoff Datetime"    -> 2016-09-02 04:51:43
# zh-cn!K3                "Correspond Handback DateTime"   -> 2016-09-02 04:52:01
# de-de!H3                "Correspond Handoff Datetime"    -> 2016-09-02 04:51:47
# de-de!K3                "Correspond Handback DateTime"   -> 2016-09-02 04:52:17

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the f8ddade8... file; column G is
# "Latest HO Xliff Generate Date".
$wsOverview.Range("G3").Value = "2016-09-02 04:51:47"

# zh-cn sheet: row 3 is the f8ddade8... file; column H is
# "Correspond Handoff Datetime", column K is "Correspond Handback DateTime".
$wsZhCn.Range("H3").Value = "2016-09-02 04:51:43"
$wsZhCn.Range("K3").Value = "2016-09-02 04:52:01"

# de-de sheet: row 3 is the f8ddade8... file; column H is
# "Correspond Handoff Datetime", column K is "Correspond Handback DateTime".
$wsDeDe.Range("H3").Value = "2016-09-02 04:51:47"
$wsDeDe.Range("K3").Value = "2016-09-02 04:52:17"
